$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.654.61"
$ws.Range("E2").Value = "  -2.45%  "

$ws.Range("D3").Value = "1.984.28"
$ws.Range("E3").Value = "  -3.79%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.10"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("E6").Value = "  -3.58%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.20"
$ws.Range("E8").Value = "  +6.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.00"
$ws.Range("E9").Value = "  +1.39%  "

$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0729"
$ws.Range("E11").Value = "  -2.73%  "

$ws.Range("E12").Value = "  -4.94%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.911"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.17"
$ws.Range("E14").Value = "  -3.47%  "

$ws.Range("D15").Value = "2.271.74"
$ws.Range("E15").Value = "  -3.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  -2.92%  "

$ws.Range("D17").Value = "1.980.98"
$ws.Range("E17").Value = "  -4.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.14"
$ws.Range("E18").Value = "  +4.91%  "

$ws.Range("D19").Value = "35.515.18"
$ws.Range("E19").Value = "  -2.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.32"
$ws.Range("E20").Value = "  -1.81%  "

$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  -2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.21"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.05"
$ws.Range("E23").Value = "  -3.68%  "

$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("E25").Value = "  -4.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +7.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.20"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.09"
$ws.Range("E28").Value = "  -4.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.46"
$ws.Range("E29").Value = "  -4.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("E30").Value = "  -2.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.13"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("E32").Value = "  -5.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0585"
$ws.Range("E33").Value = "  -1.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0902"
$ws.Range("E34").Value = "  +9.73%  "

$ws.Range("E35").Value = "  -7.25%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.80"
$ws.Range("E38").Value = "  -2.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.90"
$ws.Range("E39").Value = "  +1.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.19"
$ws.Range("E40").Value = "  -4.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.82"
$ws.Range("E41").Value = "  -3.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0209"
$ws.Range("E42").Value = "  -3.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("E43").Value = "  -5.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0888"
$ws.Range("E44").Value = "  -5.45%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.72"
$ws.Range("E45").Value = "  -3.41%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.373.79"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.41"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.46"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("E49").Value = "  +1.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").Value = "  -3.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.80"
$ws.Range("E51").Value = "  +2.39%  "
